$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3483
$ws.Range("E3").Value = 582
$ws.Range("E5").Value = 2998
$ws.Range("E6").Value = 4530
$ws.Range("E7").Value = 3387
$ws.Range("E8").Value = 10389
$ws.Range("E9").Value = 15014
$ws.Range("E10").Value = 2334
$ws.Range("E11").Value = 6497
$ws.Range("E12").Value = 3219
$ws.Range("E13").Value = 8841
$ws.Range("E14").Value = 4528
